$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K"), rows 2-24
$values = @{
    2  = 0
    3  = 2
    4  = 10
    5  = 3
    6  = 1
    7  = 5
    8  = 7
    9  = 5
    10 = 3
    11 = 7
    12 = 1
    13 = 4
    14 = 3
    15 = 2
    16 = 7
    17 = 8
    18 = 3
    19 = 5
    20 = 5
    21 = 5
    22 = 0
    23 = 1
    24 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
